# Update TrackingSheet: add BookingNumber column (G), refresh LastUpdated times.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell G1, matching the style used by the rest of row 1 (F1).
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G1").Value = "BookingNumber"

# Refresh LastUpdated timestamps.
$ws.Range("F2").Value = "2025-05-12 08:00"
$ws.Range("F3").Value = "2025-05-12 08:00"

# New (currently blank) BookingNumber cells for the data rows, carrying the
# same plain formatting as the rest of those rows (e.g. F2/F3).
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0
